# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Thu Jul 20 09:36:00 UTC 2023 with GitHub Actions".
#
# The sheet stores Price/Volume columns (D/E) as plain text, even though many
# values look numeric (e.g. "1.000", "0.08100", "103.10"). Excel normally
# auto-converts such text to numbers when set through the COM object model,
# which would silently drop meaningful trailing zeros / dot-grouping.
# Prefixing the literal with a leading apostrophe forces Excel to keep it as
# text (quote-prefixed) without altering the cells number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.327.70"
$ws.Range("E2").Value = "'  +1.17%  "
$ws.Range("D3").Value = "'1.918.62"
$ws.Range("E3").Value = "'  +0.55%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'0.8112"
$ws.Range("E5").Value = "'  +3.90%  "
$ws.Range("D6").Value = "'244.54"
$ws.Range("E6").Value = "'  +1.18%  "
$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("D8").Value = "'0.3239"
$ws.Range("E8").Value = "'  +2.69%  "
$ws.Range("D9").Value = "'27.17"
$ws.Range("E9").Value = "'  +3.95%  "
$ws.Range("D10").Value = "'0.07172"
$ws.Range("E10").Value = "'  +4.34%  "
$ws.Range("D11").Value = "'0.7856"
$ws.Range("E11").Value = "'  +6.18%  "
$ws.Range("D12").Value = "'0.08100"
$ws.Range("E12").Value = "'  +1.63%  "
$ws.Range("D13").Value = "'1.911.16"
$ws.Range("E13").Value = "'  +0.14%  "
$ws.Range("D14").Value = "'5.402"
$ws.Range("E14").Value = "'  +3.98%  "
$ws.Range("D15").Value = "'94.73"
$ws.Range("E15").Value = "'  +2.26%  "
$ws.Range("D16").Value = "'30.333.50"
$ws.Range("E16").Value = "'  +1.15%  "
$ws.Range("E17").Value = "'  +2.80%  "
$ws.Range("D18").Value = "'6.039"
$ws.Range("E18").Value = "'  +2.93%  "
$ws.Range("D19").Value = "'252.01"
$ws.Range("E19").Value = "'  +2.79%  "
$ws.Range("D20").Value = "'0.000007832"
$ws.Range("E20").Value = "'  +1.42%  "
$ws.Range("D21").Value = "'2.167.85"
$ws.Range("E21").Value = "'  +0.44%  "
$ws.Range("E22").Value = "'  +0.13%  "
$ws.Range("D23").Value = "'7.990"
$ws.Range("E23").Value = "'  +16.76%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "'  +0.12%  "
$ws.Range("D25").Value = "'0.1626"
$ws.Range("E25").Value = "'  +18.15%  "
$ws.Range("D26").Value = "'9.523"
$ws.Range("E26").Value = "'  +3.06%  "
$ws.Range("D27").Value = "'167.62"
$ws.Range("E27").Value = "'  -0.48%  "
$ws.Range("D28").Value = "'19.13"
$ws.Range("E28").Value = "'  +1.58%  "
$ws.Range("D29").Value = "'2.141"
$ws.Range("E29").Value = "'  +5.63%  "
$ws.Range("E30").Value = "'  +0.59%  "
$ws.Range("D31").Value = "'1.539"
$ws.Range("E31").Value = "'  +1.62%  "
$ws.Range("D32").Value = "'4.358"
$ws.Range("E32").Value = "'  +1.41%  "
$ws.Range("B33").Value = "'InternetComputer(DFINITY)"
$ws.Range("C33").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.139"
$ws.Range("E33").Value = "'  +1.72%  "
$ws.Range("B34").Value = "'Hedera"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.05643"
$ws.Range("E34").Value = "'  +1.69%  "
$ws.Range("D35").Value = "'1.301"
$ws.Range("E35").Value = "'  +4.10%  "
$ws.Range("D36").Value = "'0.7443"
$ws.Range("E36").Value = "'  +1.73%  "
$ws.Range("E37").Value = "'  +0.15%  "
$ws.Range("D38").Value = "'2.716"
$ws.Range("E38").Value = "'  +0.12%  "
$ws.Range("D39").Value = "'0.01953"
$ws.Range("E39").Value = "'  +1.41%  "
$ws.Range("D40").Value = "'2.820"
$ws.Range("E40").Value = "'  +1.16%  "
$ws.Range("D41").Value = "'0.4493"
$ws.Range("E41").Value = "'  +2.03%  "
$ws.Range("D42").Value = "'73.73"
$ws.Range("E42").Value = "'  +2.68%  "
$ws.Range("D43").Value = "'5.984"
$ws.Range("E43").Value = "'  -2.19%  "
$ws.Range("D44").Value = "'0.8553"
$ws.Range("E44").Value = "'  +1.79%  "
$ws.Range("D45").Value = "'1.936"
$ws.Range("E45").Value = "'  +3.63%  "
$ws.Range("E46").Value = "'  +0.12%  "
$ws.Range("D47").Value = "'1.039.92"
$ws.Range("E47").Value = "'  +5.59%  "
$ws.Range("D48").Value = "'103.10"
$ws.Range("E48").Value = "'  +2.62%  "
$ws.Range("D49").Value = "'10.04"
$ws.Range("E49").Value = "'  +3.24%  "
$ws.Range("D50").Value = "'7.653"
$ws.Range("E50").Value = "'  +2.00%  "
$ws.Range("D51").Value = "'2.074.78"
$ws.Range("E51").Value = "'  +0.89%  "
